# Update the cryptocurrency price/volume table with the latest scrape.
# (GitHub Actions scheduled refresh of cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '95.095.10'
$ws.Range("E2").Value = '  -1.36%  '

# Row 3
$ws.Range("D3").Value = '3.572.39'
$ws.Range("E3").Value = '  -0.89%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '''235.28'
$ws.Range("E5").Value = '  -2.11%  '

# Row 6
$ws.Range("D6").Value = '''658.31'
$ws.Range("E6").Value = '  +3.09%  '

# Row 7
$ws.Range("D7").Value = '''1.47'
$ws.Range("E7").Value = '  -0.86%  '

# Row 8
$ws.Range("D8").Value = '''0.398'
$ws.Range("E8").Value = '  -1.30%  '

# Row 9
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("E10").Value = '  -2.31%  '

# Row 11
$ws.Range("D11").Value = '3.567.95'
$ws.Range("E11").Value = '  -0.92%  '

# Row 12
$ws.Range("E12").Value = '  +1.01%  '

# Row 13
$ws.Range("E13").Value = '  -1.90%  '

# Row 14
$ws.Range("D14").Value = '''6.48'
$ws.Range("E14").Value = '  +0.99%  '

# Row 15
$ws.Range("D15").Value = '4.233.30'
$ws.Range("E15").Value = '  -1.35%  '

# Row 16
$ws.Range("D16").Value = '94.972.83'
$ws.Range("E16").Value = '  -1.38%  '

# Row 17
$ws.Range("E17").Value = '  -0.51%  '

# Row 18
$ws.Range("D18").Value = '3.580.56'
$ws.Range("E18").Value = '  -0.45%  '

# Row 19
$ws.Range("D19").Value = '''7.75'
$ws.Range("E19").Value = '  -7.39%  '

# Row 20
$ws.Range("D20").Value = '''12.62'
$ws.Range("E20").Value = '  -4.71%  '

# Row 21
$ws.Range("D21").Value = '''17.81'
$ws.Range("E21").Value = '  -1.69%  '

# Row 22
$ws.Range("E22").Value = '  -0.36%  '

# Row 23
$ws.Range("D23").Value = '''508.59'
$ws.Range("E23").Value = '  -1.61%  '

# Row 24
$ws.Range("D24").Value = '''0.480'
$ws.Range("E24").Value = '  -4.18%  '

# Row 25
$ws.Range("D25").Value = '''6.76'
$ws.Range("E25").Value = '  +1.15%  '

# Row 26
$ws.Range("E26").Value = '  -1.84%  '

# Row 27
$ws.Range("D27").Value = '''95.18'
$ws.Range("E27").Value = '  -1.93%  '

# Row 28
$ws.Range("D28").Value = '''12.69'
$ws.Range("E28").Value = '  +1.78%  '

# Row 29
$ws.Range("D29").Value = '3.763.01'
$ws.Range("E29").Value = '  -0.93%  '

# Row 30
$ws.Range("E30").Value = '  -1.15%  '

# Row 31
$ws.Range("D31").Value = '''11.49'
$ws.Range("E31").Value = '  -1.42%  '

# Row 32
$ws.Range("E32").Value = '  -0.44%  '

# Row 33
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  -0.04%  '

# Row 34
$ws.Range("D34").Value = '''0.999'
$ws.Range("E34").Value = '  -0.24%  '

# Row 35
$ws.Range("D35").Value = '''0.177'
$ws.Range("E35").Value = '  -2.60%  '

# Row 36
$ws.Range("D36").Value = '''31.80'
$ws.Range("E36").Value = '  +4.30%  '

# Row 37
$ws.Range("E37").Value = '  +12.57%  '

# Row 38
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Value = '''0.557'
$ws.Range("E38").Value = '  -1.99%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = '''8.49'
$ws.Range("E39").Value = '  +7.66%  '

# Row 40
$ws.Range("D40").Value = '''581.56'
$ws.Range("E40").Value = '  +1.05%  '

# Row 41
$ws.Range("E41").Value = '  +0.04%  '

# Row 42
$ws.Range("E42").Value = '  -0.91%  '

# Row 43
$ws.Range("D43").Value = '''0.906'
$ws.Range("E43").Value = '  -2.41%  '

# Row 44
$ws.Range("E44").Value = '  +3.37%  '

# Row 45
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = '''5.73'
$ws.Range("E45").Value = '  +0.86%  '

# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '''2.28'
$ws.Range("E46").Value = '  +3.54%  '

# Row 47
$ws.Range("D47").Value = '''34.28'
$ws.Range("E47").Value = '  +30.76%  '

# Row 48
$ws.Range("E48").Value = '  -1.78%  '

# Row 49
$ws.Range("E49").Value = '  -3.99%  '

# Row 50
$ws.Range("E50").Value = '  +0.82%  '

# Row 51
$ws.Range("D51").Value = '''8.17'
